$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Insert three new paragraphs at the very top of the document:
#      "Preparing standards" (bold)
#      an empty paragraph
#      an empty (bold-marked) paragraph
# ---------------------------------------------------------------------
$introXml = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Preparing standards</w:t></w:r></w:p>' +
            '<w:p ' + $wNs + '/>' +
            '<w:p ' + $wNs + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>'

$top = $d.Range(0, 0)
$null = $top.InsertXML($introXml)

# ---------------------------------------------------------------------
# 2) Wrap "microM" in the calibration-samples bullet with spell-check
#    proofErr markers, splitting the single run into three runs.
#    The paragraph's own attributes / pPr are preserved by reading them
#    back off the live document instead of being hard-coded.
# ---------------------------------------------------------------------
$calibPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Run all calibration samples*") {
        $calibPara = $p
    }
}

$owx = $calibPara.Range.WordOpenXML
$null = $owx -match '<w:body>(<w:p [^>]*>)(<w:pPr>.*?</w:pPr>)?'
$openTag = $matches[1]
$pPr = $matches[2]

$enDash = [char]0x2013
$runXml = '<w:r><w:t xml:space="preserve">Run all calibration samples (0.1 ' + $enDash + ' 100 </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>microM</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>) after initial washing</w:t></w:r>'

$newPara = $openTag + $pPr + $runXml + '</w:p>'

$null = $calibPara.Range.InsertXML($newPara)
